$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/2/2025  Through  6/8/2025"

# --- Data table updates (rows 15-30) ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("G15").Value = 2
$ws.Range("I15").Value = 13
$ws.Range("K15").Value = 62.5
$ws.Range("L15").Value = -13.333333333333
$ws.Range("N15").Value = -64.864864864864
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = -29.166666666666
$ws.Range("I16").Value = 97
$ws.Range("J16").Value = 102
$ws.Range("K16").Value = -4.901960784313
$ws.Range("L16").Value = 10.227272727272
$ws.Range("M16").Value = -21.774193548387
$ws.Range("N16").Value = -73.497267759562
$ws.Range("C17").Value = 12
$ws.Range("E17").Value = 140
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 15.625
$ws.Range("I17").Value = 205
$ws.Range("J17").Value = 216
$ws.Range("K17").Value = -5.092592592592
$ws.Range("L17").Value = 1.990049751243
$ws.Range("M17").Value = 65.322580645161
$ws.Range("N17").Value = -48.621553884711
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 42.857142857142
$ws.Range("I18").Value = 45
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = -25
$ws.Range("L18").Value = -44.444444444444
$ws.Range("M18").Value = 4.651162790697
$ws.Range("N18").Value = -89.311163895486
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -32.432432432432
$ws.Range("I19").Value = 167
$ws.Range("J19").Value = 156
$ws.Range("K19").Value = 7.051282051282
$ws.Range("L19").Value = -1.764705882352
$ws.Range("M19").Value = 26.515151515151
$ws.Range("N19").Value = 4.375
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 57.142857142857
$ws.Range("I20").Value = 28
$ws.Range("J20").Value = 42
$ws.Range("L20").Value = -20
$ws.Range("M20").Value = 21.739130434782
$ws.Range("N20").Value = -77.419354838709
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 8.695652173913
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = -9.090909090909
$ws.Range("I21").Value = 559
$ws.Range("J21").Value = 588
$ws.Range("K21").Value = -4.931972789115
$ws.Range("L21").Value = -5.733558178752
$ws.Range("M21").Value = 21.258134490238
$ws.Range("N21").Value = -63.416230366492
$ws.Range("C14").Copy($ws.Range("G22"))
$ws.Range("E14").Copy($ws.Range("H22"))
$ws.Range("L22").Value = -42.857142857142
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = -11.764705882352
$ws.Range("I23").Value = 102
$ws.Range("J23").Value = 117
$ws.Range("K23").Value = -12.820512820512
$ws.Range("L23").Value = 5.154639175257
$ws.Range("M23").Value = 45.714285714285
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -28.571428571428
$ws.Range("F24").Value = 79
$ws.Range("G24").Value = 85
$ws.Range("H24").Value = -7.058823529411
$ws.Range("I24").Value = 416
$ws.Range("J24").Value = 385
$ws.Range("K24").Value = 8.051948051948
$ws.Range("L24").Value = 19.884726224783
$ws.Range("M24").Value = 46.996466431095
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 80
$ws.Range("I25").Value = 86
$ws.Range("J25").Value = 75
$ws.Range("K25").Value = 14.666666666666
$ws.Range("L25").Value = 10.25641025641
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 23
$ws.Range("E26").Value = -21.739130434782
$ws.Range("F26").Value = 62
$ws.Range("G26").Value = 69
$ws.Range("H26").Value = -10.144927536231
$ws.Range("I26").Value = 285
$ws.Range("J26").Value = 295
$ws.Range("K26").Value = -3.389830508474
$ws.Range("L26").Value = 6.343283582089
$ws.Range("M26").Value = -24
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 15
$ws.Range("K27").Value = 36.363636363636
$ws.Range("L27").Value = -31.818181818181
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 37.5
$ws.Range("I28").Value = 36
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = 44
$ws.Range("L28").Value = 56.521739130434
$ws.Range("D16").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 2
$ws.Range("E16").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 11
$ws.Range("K29").Value = -36.363636363636
$ws.Range("L29").Value = -36.363636363636
$ws.Range("M29").Value = -53.333333333333
$ws.Range("N29").Value = -85.106382978723
$ws.Range("D16").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 2
$ws.Range("E16").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 8
$ws.Range("K30").Value = -25
$ws.Range("L30").Value = -33.333333333333
$ws.Range("M30").Value = -45.454545454545
$ws.Range("N30").Value = -86.046511627907
